$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 355; this shifts the existing rows 355-430 down to 356-431
$ws.Rows.Item(355).Insert()

# Populate the new row 355 with its data (same non-changed fields as the former
# row 355, which is now row 356, plus the new values for this record)
$ws.Cells.Item(355, 1).Value = 10
$ws.Cells.Item(355, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(355, 3).Value = "La Araucanía"
$ws.Cells.Item(355, 4).Value = 44995
$ws.Cells.Item(355, 5).Value = 9
$ws.Cells.Item(355, 6).Value = 100112017
$ws.Cells.Item(355, 7).Value = "Apio"
$ws.Cells.Item(355, 8).Value = "Americana (o)"
$ws.Cells.Item(355, 9).Value = "Primera"
$ws.Cells.Item(355, 10).Value = 65
$ws.Cells.Item(355, 11).Value = 13000
$ws.Cells.Item(355, 12).Value = 13000
$ws.Cells.Item(355, 13).Value = 13000
$ws.Cells.Item(355, 14).Value = "$/docena de matas"
$ws.Cells.Item(355, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(355, 16).Value = 2167
$ws.Cells.Item(355, 17).Value = 6
$ws.Cells.Item(355, 18).Value = "Hortaliza"

# Match the existing date-column number format used by the rest of column D
$ws.Cells.Item(355, 4).NumberFormat = $ws.Cells.Item(356, 4).NumberFormat
